# emis_preds/inputs/inputs.xlsx  --  "Add mean calc, add new summ stats to plot, A, B, C"
#
# The "Slurry" sheet's digestate rows are being re-organised:
#   - the old literature-derived rows for cattle/pig/mix digestate (rows 6-8,
#     labelled "Digestate CL/PL/ML") are dropped
#   - the old 10th/90th percentile rows (rows 9-10, "Digestate low/high") move
#     up to rows 6-7 and are relabelled "Digestate B low/high"
#   - a brand new "Digestate C" block (rows 8-10: mean, low, high) is added
#   - an (as yet empty, "0.00"-formatted) row 11 is prepared for a mean calc
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# ---------------------------------------------------------------------------
# Row 6: was "digestate cattle lit" / "Digestate CL" / 4.03 / 7.33
#        -> "digestate cris 10th" (rich) / "Digestate B low" / 5.45 / 7.864
$c = $ws.Range("A6")
$c.Value = "digestate cris 10th"
$c.Characters(18, 2).Font.Superscript = $true
$ws.Range("B6").Value = "Digestate B low"
$ws.Range("C6").NumberFormat = "General"
$ws.Range("C6").Value = 5.45
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Value = 7.864

# ---------------------------------------------------------------------------
# Row 7: was "digestate pig lit" / "Digestate PL" / 3.68 / 7.67
#        -> "digestate cris 90th" (rich) / "Digestate B high" / 9.114 / 8.282
$c = $ws.Range("A7")
$c.Value = "digestate cris 90th"
$c.Characters(18, 2).Font.Superscript = $true
$ws.Range("B7").Value = "Digestate B high"
$ws.Range("C7").NumberFormat = "General"
$ws.Range("C7").Value = 9.114
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Value = 8.282

# ---------------------------------------------------------------------------
# Row 8: was "digestate mix lit" / "Digestate ML" / 4.55 / 8.16
#        -> "digestate lit" (plain) / "Digestate C" / 3.87 / 7.97
$ws.Range("A8").Value = "digestate lit"
$ws.Range("B8").Value = "Digestate C"
$ws.Range("C8").NumberFormat = "General"
$ws.Range("C8").Value = 3.87
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Value = 7.97

# ---------------------------------------------------------------------------
# Row 9: was "digestate cris 10th" / "Digestate low" / 5.45 / 7.864
#        -> "digestate lit 10th" (rich) / "Digestate C low" / 1.99 / 7.41
$c = $ws.Range("A9")
$c.Value = "digestate lit 10th"
$c.Characters(1, 16).Font.Superscript = $false
$c.Characters(17, 2).Font.Superscript = $true
$ws.Range("B9").Value = "Digestate C low"
$ws.Range("C9").NumberFormat = "General"
$ws.Range("C9").Value = 1.99
$ws.Range("D9").NumberFormat = "0.00"
$ws.Range("D9").Value = 7.41

# ---------------------------------------------------------------------------
# Row 10: was "digestate cris 90th" / "Digestate high" / 9.114 / 8.282
#         -> "digestate lit 90th" (rich) / "Digestate C high" / 5.54 / 8.53
$c = $ws.Range("A10")
$c.Value = "digestate lit 90th"
$c.Characters(1, 16).Font.Superscript = $false
$c.Characters(17, 2).Font.Superscript = $true
$ws.Range("B10").Value = "Digestate C high"
$ws.Range("C10").NumberFormat = "General"
$ws.Range("C10").Value = 5.54
$ws.Range("D10").NumberFormat = "0.00"
$ws.Range("D10").Value = 8.53

# ---------------------------------------------------------------------------
# Row 11: new, blank placeholder cells for the mean-calc row, "0.00" formatted
$ws.Range("C11").NumberFormat = "0.00"
$ws.Range("D11").NumberFormat = "0.00"

# Put the selection where the author left it
$ws.Range("B11").Select()
